$d = $word.ActiveDocument

$pairs = @(
    @("38×58=2204", "26×14=364"),
    @("27×56=1512", "53×43=2279"),
    @("70×99=6930", "25×34=850"),
    @("84×47=3948", "21×42=882"),
    @("45×75=3375", "74×97=7178"),
    @("28×17=476", "83×87=7221"),
    @("89×48=4272", "85×12=1020"),
    @("77×30=2310", "27×75=2025"),
    @("77×68=5236", "15×87=1305"),
    @("27×65=1755", "22×97=2134"),
    @("52×32=1664", "34×40=1360"),
    @("89×35=3115", "27×19=513"),
    @("96×23=2208", "14×32=448"),
    @("27×86=2322", "92×42=3864"),
    @("11×63=693", "29×50=1450"),
    @("62×68=4216", "41×66=2706"),
    @("66×47=3102", "45×21=945"),
    @("54×96=5184", "65×91=5915"),
    @("80×11=880", "97×32=3104"),
    @("50×50=2500", "76×35=2660"),
    @("60×62=3720", "84×21=1764"),
    @("98×52=5096", "53×27=1431"),
    @("88×85=7480", "75×23=1725"),
    @("81×89=7209", "94×85=7990"),
    @("18×75=1350", "60×13=780")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
